$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Update existing row 60 (was Web/Student) -> (Android/Teacher) for the DoubtForum Tabs test case
$ws.Range("B60").Value = "Android"
$ws.Range("D60").Value = "Teacher"

# Insert a new row at position 61 (pushes the old row 61 down to row 62),
# then populate it with a new "Web / Student" doubt-forum test case.
$ws.Rows("61").Insert()
$ws.Range("B61").Value = "Web"
$ws.Range("C61").Value = "verifyCountOfDoubtTabs"
$ws.Range("D61").Value = "Student"
$ws.Range("E61").Value = "DoubtForum Tabs"
$ws.Range("F61").Value = "Doubts, My Doubts, Post For Me"

# Match the workbook's recorded active-cell selection after the edit
$ws.Range("C60").Select()
